$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the first sheet: Plan1 -> trajetoria
# ---------------------------------------------------------------------------
$traj = $wb.Worksheets.Item("Plan1")
$traj.Name = "trajetoria"

# ---------------------------------------------------------------------------
# 2) Fix up the line chart on "trajetoria" so its series formulas point at
#    the renamed sheet instead of the old "Plan1" name.
# ---------------------------------------------------------------------------
$chartObj = $traj.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(trajetoria!`$B`$1,trajetoria!`$A`$2:`$A`$18,trajetoria!`$B`$2:`$B`$18,1)"

# ---------------------------------------------------------------------------
# 3) Create the new "skills" and "formação" sheets (while the old "Plan2"
#    sheet still exists) so the internal sheetId counter lands on the same
#    values as the authored workbook (skills=3, formação=4), then drop the
#    old "Plan2" ("Contato") sheet entirely - its content is being replaced.
# ---------------------------------------------------------------------------
$skills = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $traj)
$skills.Name = "skills_tmp"

$formacao = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $skills)
$formacao.Name = "formação"

[void]$wb.Worksheets.Item("Plan2").Delete()
$skills.Name = "skills"

# ---------------------------------------------------------------------------
# 4) Populate "skills" with the Skill / Grau table.
# ---------------------------------------------------------------------------
$skills.Range("A1").Value = "Skill"
$skills.Range("A1").Font.Bold = $true

$skills.Range("B1").Value = "Grau"
$skills.Range("B1").Font.Bold = $true
$skills.Range("B1").HorizontalAlignment = -4108

$skillRows = @(
    @("LaTeX", 0.8),
    @("R", 0.7),
    @("RMarkDown", 0.7),
    @("Sparklyr", 0.7),
    @("Dplyr e tidyverse", 0.7),
    @("APIs", 0.6),
    @("Excel e VBA", 0.7),
    @("Econometria e Estatística", 0.7),
    @("Machine Learning", 0.6),
    @("Bancos de Dados", 0.6),
    @("Git e Github", 0.5),
    @("SQL", 0.5),
    @("Google Colab", 0.55000000000000004),
    @("Google Cloud Platform", 0.3),
    @("AWS", 0.3),
    @("Python", 0.1),
    @("JIRA", 0.3),
    @("Confluence", 0.8),
    @("ETL", 0.8),
    @("PmBok", 0.8),
    @("Finanças e pricing", 0.7),
    @("Demand forecast", 0.7)
)

$r = 2
foreach ($row in $skillRows) {
    $skills.Cells.Item($r, 1).Value = $row[0]
    $skills.Cells.Item($r, 2).Value = $row[1]
    $skills.Cells.Item($r, 2).NumberFormat = "0%"
    $skills.Cells.Item($r, 2).HorizontalAlignment = -4108
    $r++
}

$skills.Columns.Item(1).ColumnWidth = 21.81640625
$skills.Columns.Item(2).ColumnWidth = 8.7265625

# ---------------------------------------------------------------------------
# 5) Restore view state: "skills" is the selected/active tab (as "Plan2" was
#    in the source workbook), zoomed to 85%, cursor resting below the table.
# ---------------------------------------------------------------------------
[void]$skills.Activate()
$excel.ActiveWindow.Zoom = 85
[void]$skills.Range("A25").Select()

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
